# Refresh the cryptos price/volume table (and two row re-orderings) to
# match the latest coinranking.com snapshot.
# Commit message: "Updated cryptos list on Wed Dec  6 11:49:50 UTC 2023
# with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    # D-column "prices" are plain text in this sheet (some look like
    # numbers, e.g. "6.00" / "1.00" / "229.69"); a leading apostrophe
    # forces Excel to store the literal text instead of silently
    # coercing it to a Double and losing trailing zeros / formatting.
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

# Row 2
Set-TextValue 2 4 '43.719.09'
$ws.Cells.Item(2, 5).Value = '  +4.83%  '

# Row 3
Set-TextValue 3 4 '2.245.73'
$ws.Cells.Item(3, 5).Value = '  +2.02%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.24%  '

# Row 5
Set-TextValue 5 4 '229.69'
$ws.Cells.Item(5, 5).Value = '  +0.08%  '

# Row 6
Set-TextValue 6 4 '0.628'
$ws.Cells.Item(6, 5).Value = '  +2.16%  '

# Row 7
Set-TextValue 7 4 '63.27'
$ws.Cells.Item(7, 5).Value = '  +4.84%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.06%  '

# Row 9
Set-TextValue 9 4 '0.433'
$ws.Cells.Item(9, 5).Value = '  +8.42%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +12.97%  '

# Row 11
Set-TextValue 11 4 '56.25'
$ws.Cells.Item(11, 5).Value = '  -1.07%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +2.42%  '

# Row 13
Set-TextValue 13 4 '25.66'
$ws.Cells.Item(13, 5).Value = '  +16.68%  '

# Row 14
Set-TextValue 14 4 '2.581.65'
$ws.Cells.Item(14, 5).Value = '  +2.13%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  +1.65%  '

# Row 16
Set-TextValue 16 4 '5.93'
$ws.Cells.Item(16, 5).Value = '  +6.89%  '

# Row 17
Set-TextValue 17 4 '0.818'
$ws.Cells.Item(17, 5).Value = '  +3.21%  '

# Row 18
Set-TextValue 18 4 '2.238.77'
$ws.Cells.Item(18, 5).Value = '  +1.44%  '

# Row 19
Set-TextValue 19 4 '43.556.10'
$ws.Cells.Item(19, 5).Value = '  +4.67%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +11.73%  '

# Row 21
Set-TextValue 21 4 '73.34'
$ws.Cells.Item(21, 5).Value = '  +2.14%  '

# Row 22
Set-TextValue 22 4 '6.01'
$ws.Cells.Item(22, 5).Value = '  -0.14%  '

# Row 23
Set-TextValue 23 4 '252.05'
$ws.Cells.Item(23, 5).Value = '  +4.72%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +0.10%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +3.74%  '

# Row 26
Set-TextValue 26 4 '2.33'
$ws.Cells.Item(26, 5).Value = '  +2.66%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +3.20%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +2.00%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'Kaspa'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 29 4 '0.137'
$ws.Cells.Item(29, 5).Value = '  -1.42%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'EthereumClassic'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 30 4 '20.62'
$ws.Cells.Item(30, 5).Value = '  +4.76%  '

# Row 31
Set-TextValue 31 4 '2.84'
$ws.Cells.Item(31, 5).Value = '  +9.18%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -3.75%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +2.38%  '

# Row 34
Set-TextValue 34 4 '0.0673'
$ws.Cells.Item(34, 5).Value = '  +4.65%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +1.89%  '

# Row 36
Set-TextValue 36 4 '4.92'
$ws.Cells.Item(36, 5).Value = '  -0.44%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +9.09%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +5.55%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -1.04%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +4.24%  '

# Row 41
Set-TextValue 41 4 '1.00'
$ws.Cells.Item(41, 5).Value = '  +0.09%  '

# Row 42
Set-TextValue 42 4 '17.48'
$ws.Cells.Item(42, 5).Value = '  +9.18%  '

# Row 43
Set-TextValue 43 4 '8.20'
$ws.Cells.Item(43, 5).Value = '  -4.52%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +0.61%  '

# Row 45
Set-TextValue 45 4 '96.66'
$ws.Cells.Item(45, 5).Value = '  +0.29%  '

# Row 46
Set-TextValue 46 4 '0.000210'
$ws.Cells.Item(46, 5).Value = '  -13.86%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'FTXToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 47 4 '4.34'
$ws.Cells.Item(47, 5).Value = '  -0.06%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 48 4 '1.18'
$ws.Cells.Item(48, 5).Value = '  -1.65%  '

# Row 49
Set-TextValue 49 4 '1.446.18'
$ws.Cells.Item(49, 5).Value = '  -0.22%  '

# Row 50
Set-TextValue 50 4 '2.29'
$ws.Cells.Item(50, 5).Value = '  +4.35%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -0.51%  '
